$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "51.768.77"
$ws.Range("E2").Value = "  +0.10%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.804.58"
$ws.Range("E3").Value = "  +0.88%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "354.79"
$ws.Range("E5").Value = "  -0.42%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "111.53"
$ws.Range("E6").Value = "  +2.43%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.558"
$ws.Range("E7").Value = "  +0.79%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.633"
$ws.Range("E9").Value = "  +8.00%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "40.39"
$ws.Range("E10").Value = "  +1.93%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -1.45%  "

# Row 12 - Dogecoin
Set-TextValue $ws.Range("D12") "0.0838"

# Row 13 - Chainlink
$ws.Range("E13").Value = "  +2.74%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +2.47%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "3.245.89"
$ws.Range("E15").Value = "  +0.82%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "2.799.73"
$ws.Range("E16").Value = "  +0.11%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.947"
$ws.Range("E17").Value = "  +1.49%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "51.788.46"
$ws.Range("E18").Value = "  +0.20%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  +1.71%  "

# Row 20 - ImmutableX
Set-TextValue $ws.Range("D20") "3.18"
$ws.Range("E20").Value = "  +3.06%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D21") "13.66"
$ws.Range("E21").Value = "  +3.47%  "

# Row 22 - ShibaInu
$ws.Range("E22").Value = "  +0.97%  "

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") "70.59"
$ws.Range("E23").Value = "  +0.60%  "

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") "268.83"
$ws.Range("E24").Value = "  +0.61%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +1.32%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.15%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "26.16"

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -2.71%  "

# Row 29 - InjectiveProtocol
Set-TextValue $ws.Range("D29") "38.76"
$ws.Range("E29").Value = "  +11.41%  "

# Row 30 - Cosmos
Set-TextValue $ws.Range("D30") "10.39"
$ws.Range("E30").Value = "  +1.98%  "

# Row 31 - Toncoin
$ws.Range("E31").Value = "  +3.72%  "

# Row 32 - OKB
Set-TextValue $ws.Range("D32") "52.34"
$ws.Range("E32").Value = "  +1.33%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "6.14"
$ws.Range("E33").Value = "  +0.27%  "

# Row 34 - RenderToken
Set-TextValue $ws.Range("D34") "5.64"
$ws.Range("E34").Value = "  +8.89%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  +5.93%  "

# Row 36 - VeChain
$ws.Range("E36").Value = "  -0.36%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  -0.09%  "

# Row 38 - Celestia
Set-TextValue $ws.Range("D38") "18.86"
$ws.Range("E38").Value = "  +0.09%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  +2.60%  "

# Row 40 - LidoDAOToken
Set-TextValue $ws.Range("D40") "3.16"
$ws.Range("E40").Value = "  +0.91%  "

# Row 41 - Stellar
$ws.Range("E41").Value = "  +1.24%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -1.94%  "

# Row 43 - was WEMIXToken, now EnergySwap
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D43") "22.12"
$ws.Range("E43").Value = "  +2.00%  "

# Row 44 - was Monero, now WEMIXToken
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D44") "2.21"
$ws.Range("E44").Value = "  -0.09%  "

# Row 45 - was EnergySwap, now Monero
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D45") "119.82"
$ws.Range("E45").Value = "  -0.14%  "

# Row 46 - NEARProtocol
Set-TextValue $ws.Range("D46") "3.42"
$ws.Range("E46").Value = "  +4.79%  "

# Row 47 - Maker
Set-TextValue $ws.Range("D47") "2.115.60"
$ws.Range("E47").Value = "  +1.20%  "

# Row 48 - ApeXProtocol
Set-TextValue $ws.Range("D48") "2.42"
$ws.Range("E48").Value = "  +6.53%  "

# Row 49 - SEI
Set-TextValue $ws.Range("D49") "0.953"
$ws.Range("E49").Value = "  +0.04%  "

# Row 50 - was THORChain, now TrustWalletToken
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D50") "1.38"
$ws.Range("E50").Value = "  +7.97%  "

# Row 51 - was TrustWalletToken, now THORChain
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D51") "5.48"
$ws.Range("E51").Value = "  -1.16%  "
